$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("max-arrecad")

# Re-assign the category labels (column A) among the rows that are tied
# on the "max arrecadado" value, fixing the ranking so it matches the
# author/mention database for the 2024 analysis year. The numeric values
# in column B are untouched.

$ws.Range("A2").Value = "jogos"
$ws.Range("A3").Value = "humor"
$ws.Range("A4").Value = "religiosidade"
$ws.Range("A5").Value = "terror"
$ws.Range("A6").Value = "hqmix"
$ws.Range("A7").Value = "angelo_agostini"

$ws.Range("A17").Value = "zine"
$ws.Range("A18").Value = "herois"
